# Update gh-pages output values (regenerated scrape numbers) across the
# four worksheets of the workbook: 展览, 演出, 本地生活, 全部类型.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 184
$ws1.Range("F6").Value  = 8
$ws1.Range("F7").Value  = 758
$ws1.Range("F8").Value  = 59
$ws1.Range("F9").Value  = 9842
$ws1.Range("F13").Value = 2397
$ws1.Range("F14").Value = 2670
$ws1.Range("F17").Value = 2097
$ws1.Range("F22").Value = 74
$ws1.Range("F25").Value = 160
$ws1.Range("F26").Value = 602
$ws1.Range("F32").Value = 1687
$ws1.Range("F33").Value = 2834
$ws1.Range("F35").Value = 1000
$ws1.Range("F36").Value = 362
$ws1.Range("C39").Value = "北京·第一届世界计划pjsk only同人展"
$ws1.Range("F39").Value = 57
$ws1.Range("F43").Value = 29

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 12

# ---- Sheet 3: 本地生活 ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 726
$ws3.Range("F3").Value = 958
$ws3.Range("F5").Value = 1824

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 726
$ws4.Range("F3").Value  = 958
$ws4.Range("F9").Value  = 184
$ws4.Range("F10").Value = 758
$ws4.Range("F11").Value = 59
$ws4.Range("F12").Value = 9842
$ws4.Range("F18").Value = 2397
$ws4.Range("F19").Value = 2671
$ws4.Range("F21").Value = 2097
$ws4.Range("F28").Value = 160
$ws4.Range("F29").Value = 602
$ws4.Range("F34").Value = 1687
$ws4.Range("F36").Value = 2834
$ws4.Range("F37").Value = 1000
$ws4.Range("F39").Value = 362
$ws4.Range("C45").Value = "北京·第一届世界计划pjsk only同人展"
$ws4.Range("F45").Value = 57
$ws4.Range("F48").Value = 29
